# Update FAST_holdings.xlsx:
#   - bump the "as of" date in the confidentiality footer from 2021-04-08 to 2021-04-09
#   - refresh the Weight / Percent Change figures in D2:E10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cell writes actually stick.
$ws.Unprotect("D382")

# --- Footer text (shared string referenced by A13) ---------------------
$footer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."
$ws.Range("A13").Value = $footer

# --- Weight (D) / Percent Change (E) values -----------------------------
$ws.Range("D2").Value = 0.1031466095110457
$ws.Range("E2").Value = -0.01059990978800196

$ws.Range("D3").Value = 0.1081472610947184
$ws.Range("E3").Value = 0.006498493531045035

$ws.Range("D4").Value = 0.116769514793387
$ws.Range("E4").Value = 0.007916833266693457

$ws.Range("D5").Value = 0.1364633175822884
$ws.Range("E5").Value = 0.004897959183673306

$ws.Range("D6").Value = 0.1327994495177891
$ws.Range("E6").Value = 0.003972194637537285

$ws.Range("D7").Value = 0.1436924331624459
$ws.Range("E7").Value = 0.005554491476728485

$ws.Range("D8").Value = 0.1298039732128766
$ws.Range("E8").Value = 0.002902757619738461

$ws.Range("D9").Value = 0.129177441125449
$ws.Range("E9").Value = -0.0001290509912729165

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.002888048711710001
